$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A second "Docente responsável" (responsible teacher) was added to the course
# sheet. That pushed the whole "Programa resumido / Short syllabus / ..." block
# down by one row. Insert a fresh row right after the existing teacher row (13)
# and fill it in with the new teacher's name in columns B and C (matching the
# layout used for the existing teacher row).
$ws.Rows.Item(14).Insert()

$ws.Range("B14").Value = "5817692 - Katia Cristiane Gandolpho Candioto"
$ws.Range("C14").Value = "5817692 - Katia Cristiane Gandolpho Candioto"
